# Update "想去人数" (interest count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 422   # 合肥·W·A第五人格同人only2.0: 421 -> 422
$wsExpo.Range("F3").Value = 2754  # 合肥·第九届环形宇宙动漫游戏嘉年华: 2740 -> 2754
$wsExpo.Range("F4").Value = 125   # 合肥·MAX特摄同人only2.0: 122 -> 125

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 422    # 合肥·W·A第五人格同人only2.0: 421 -> 422
$wsAll.Range("F7").Value = 2754   # 合肥·第九届环形宇宙动漫游戏嘉年华: 2740 -> 2754
$wsAll.Range("F8").Value = 125    # 合肥·MAX特摄同人only2.0: 122 -> 125
